# [PHOENIX-5860] Updated Create Trade License Screen
# Adds a new "tradeDetails" worksheet (tradeLicenseTestData.xlsx) with the
# trade-details test fixture row, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# --- tradeLocationDetails: the sheet was previously the active/selected
#     sheet (selection at C2). Once the new sheet becomes active this one
#     keeps its own remembered selection, which moved to A2.
$wsLocation = $wb.Worksheets.Item("tradeLocationDetails")
$wsLocation.Activate()
[void]$wsLocation.Range("A2").Select()

# --- Add the new "tradeDetails" sheet as the last (3rd) tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "tradeDetails"

# Column widths matching the source fixture's layout
$ws.Columns.Item(1).ColumnWidth = 25.56
$ws.Columns.Item(2).ColumnWidth = 17.78
$ws.Columns.Item(3).ColumnWidth = 11.52
$ws.Columns.Item(4).ColumnWidth = 19.04
$ws.Columns.Item(5).ColumnWidth = 18.89
$ws.Columns.Item(6).ColumnWidth = 27.09
$ws.Columns.Item(7).ColumnWidth = 17.22
$ws.Columns.Item(8).ColumnWidth = 27.78

# Header row (row 1) - column titles for the trade-details test data
$ws.Range("A1").Value = "dataName"
$ws.Range("B1").Value = "tradeTitle"
$ws.Range("C1").Value = "tradeType"
$ws.Range("D1").Value = "tradeCategory"
$ws.Range("E1").Value = "tradeSubCategory"
$ws.Range("F1").Value = "tradeAreaWeightOfPremises"
$ws.Range("G1").Value = "remarks"
$ws.Range("H1").Value = "tradeCommencementDate"

# Data row (row 2) - sample trade-details test fixture
$ws.Range("A2").Value = "tradeDetailsTradeLicense"
$ws.Range("B2").Value = "New Trade"
$ws.Range("C2").Value = "Permanent"
$ws.Range("D2").Value = "Flammables"
$ws.Range("E2").Value = "Petrol or Diesel bunks"
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = "Trade Data Filled"

# tradeCommencementDate must stay literal text ("11/12/2016"), not be
# reinterpreted as a date serial - format the cell as Text before writing it.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "11/12/2016"

# The rest of the sheet (including the header's H1) uses the same Text
# number format as the source fixture; F2's numeric value was already
# written above so applying "@" afterwards keeps it a real number.
$ws.Range("A1:G2").NumberFormat = "@"
$ws.Range("H1").NumberFormat = "@"

# This new sheet becomes the active tab, with E6 as its remembered selection
[void]$ws.Range("E6").Select()
